$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 414.57144
$ws.Range("I9").Value = 414.57144
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 414.57144
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -245.57144
$ws.Range("N9").Value = $null

$ws.Range("H21").Value = 7500
$ws.Range("I21").Value = 7500
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 7500
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -7032
$ws.Range("N21").Value = $null

$ws.Range("H23").Value = 7500
$ws.Range("I23").Value = 7500
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 7500
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -7266
$ws.Range("N23").Value = $null

$ws.Range("H132").Value = 3097.8333
$ws.Range("I132").Value = 3106.7273
$ws.Range("K132").Value = 9320.1819
$ws.Range("M132").Value = -6790.1819

$ws.Range("H138").Value = 3067.7646
$ws.Range("J138").Value = 3867.6667
$ws.Range("L138").Value = 11603.0001
$ws.Range("N138").Value = -21883.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4532950
$ws.Range("I32").Value = 4378759.5
$ws.Range("K32").Value = 4378759.5
$ws.Range("M32").Value = -4378472.5

$ws.Range("H46").Value = 4847.25
$ws.Range("J46").Value = 4847.25
$ws.Range("L46").Value = 4847.25
$ws.Range("N46").Value = -5485.25

$ws.Range("H74").Value = 2114.2856
$ws.Range("I74").Value = 2926.6667
$ws.Range("J74").Value = 1505
$ws.Range("K74").Value = 2926.6667
$ws.Range("L74").Value = 1505
$ws.Range("M74").Value = -2052.6667
$ws.Range("N74").Value = -3253

$ws.Range("H77").Value = 2114.2856
$ws.Range("I77").Value = 2926.6667
$ws.Range("J77").Value = 1505
$ws.Range("K77").Value = 14633.3335
$ws.Range("L77").Value = 7525
$ws.Range("M77").Value = -10265.3335
$ws.Range("N77").Value = -16261

$ws.Range("H102").Value = 499.66666
$ws.Range("I102").Value = 499.66666
$ws.Range("K102").Value = 499.66666
$ws.Range("M102").Value = 1122.33334

$ws.Range("H124").Value = 90000
$ws.Range("J124").Value = 90000
$ws.Range("L124").Value = 90000
$ws.Range("N124").Value = -99820

$ws.Range("H132").Value = 1298.75
$ws.Range("I132").Value = 1331.6666
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 3994.9998
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -1464.9998
$ws.Range("N132").Value = -8660

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 24996.334
$ws.Range("J33").Value = 24996.334
$ws.Range("L33").Value = 24996.334
$ws.Range("N33").Value = -25668.334

$ws.Range("H94").Value = 4416
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = $null

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = $null
$ws.Range("N99").Value = $null

$ws.Range("H105").Value = 2800
$ws.Range("I105").Value = 3000
$ws.Range("K105").Value = 3000
$ws.Range("M105").Value = -1253

$ws.Range("H107").Value = 1524.25
$ws.Range("I107").Value = 1550.8334
$ws.Range("K107").Value = 1550.8334
$ws.Range("M107").Value = 369.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3589.7273
$ws.Range("I105").Value = 2494
$ws.Range("J105").Value = 4215.857
$ws.Range("K105").Value = 2494
$ws.Range("L105").Value = 4215.857
$ws.Range("M105").Value = -747
$ws.Range("N105").Value = -7709.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 5000
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = $null

$ws.Range("H78").Value = 5000
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = $null

$ws.Range("H137").Value = 2489.4736
$ws.Range("I137").Value = 1138.8
$ws.Range("J137").Value = 2971.8572
$ws.Range("K137").Value = 3416.4
$ws.Range("L137").Value = 8915.571599999999
$ws.Range("M137").Value = 1683.6
$ws.Range("N137").Value = -19115.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 998.625
$ws.Range("I113").Value = 998.625
$ws.Range("K113").Value = 998.625
$ws.Range("M113").Value = 1171.375

$ws.Range("H132").Value = 8091.3335
$ws.Range("I132").Value = 8244.875
$ws.Range("K132").Value = 24734.625
$ws.Range("M132").Value = -22204.625

$ws.Range("H134").Value = 37887.5
$ws.Range("J134").Value = 37887.5
$ws.Range("L134").Value = 113662.5
$ws.Range("N134").Value = -118732.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2793.75
$ws.Range("I22").Value = 3640
$ws.Range("J22").Value = 1383.3334
$ws.Range("K22").Value = 3640
$ws.Range("L22").Value = 1383.3334
$ws.Range("M22").Value = -3345
$ws.Range("N22").Value = -1973.3334

$ws.Range("H27").Value = 2793.75
$ws.Range("I27").Value = 3640
$ws.Range("J27").Value = 1383.3334
$ws.Range("K27").Value = 3640
$ws.Range("L27").Value = 1383.3334
$ws.Range("M27").Value = -3533
$ws.Range("N27").Value = -1597.3334

$ws.Range("H132").Value = 4166.3335
$ws.Range("I132").Value = 4166.3335
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12499.0005
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9969.000499999998
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 30000
$ws.Range("I75").Value = 30000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 30000
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -29064
$ws.Range("N75").Value = $null

$ws.Range("H78").Value = 30000
$ws.Range("I78").Value = 30000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 90000
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -85320
$ws.Range("N78").Value = $null

$ws.Range("H122").Value = 1271.4286
$ws.Range("I122").Value = 816.6667
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 2450.0001
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -0.0001000000002022716
$ws.Range("N122").Value = -16900

$ws.Range("H132").Value = 5600
$ws.Range("I132").Value = 5600
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 16800
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -14270
$ws.Range("N132").Value = $null
